# Adding tutorial part 3/4 row to the Machine Learning tracker sheet.
#
# Summary of the change:
#  - D5 (the "Part 3" row) gets a Completed date (43080 = 2017-12-11).
#  - C5's URL text becomes a real clickable hyperlink (self-referencing,
#    same pattern as the existing C2 hyperlink).
#  - A brand-new row is inserted at row 6 for "Machine Learning Crash
#    Course: Part 4" (pushing the old row 6, "Python Notebook" example,
#    down to row 7). The new C6 cell is styled like a link (shrink-to-fit
#    "Hyperlink" cell style) but, like C3/C4, is NOT an actual hyperlink.
#    D6 is blank but carries the date number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 5: give it a Completed date, matching the format used elsewhere
# in column D (copy the format from D4 first since D5 has no style yet). ---
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial($xlPasteFormats)
$ws.Range("D5").Value = 43080

# --- Insert the brand new row 6 (old row 6 shifts down to row 7). ---
$ws.Rows("6:6").Insert()

$ws.Range("A6").Value = "Machine Learning"
$ws.Range("B6").Value = "Machine Learning Crash Course: Part 4"
$ws.Range("C6").Value = "https://ml.berkeley.edu/blog/2017/07/13/tutorial-4/"

# Match C6/D6 formatting to the rest of the table (copy from the
# still-untouched row 4, which already has the right styles).
$ws.Range("C4").Copy()
$ws.Range("C6").PasteSpecial($xlPasteFormats)
$ws.Range("D4").Copy()
$ws.Range("D6").PasteSpecial($xlPasteFormats)
$ws.Range("D6").ClearContents()

# --- Turn C5's URL text into a real hyperlink pointing at itself. ---
$ws.Hyperlinks.Add($ws.Range("C5"), "https://ml.berkeley.edu/blog/2017/02/04/tutorial-3/")
# Adding the hyperlink nudges the cell style; restore it to match the
# existing Hyperlink-styled cells (e.g. C4) exactly.
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- Leave the selection where the author's last save left it. ---
$ws.Activate() | Out-Null
$ws.Range("E17").Select() | Out-Null
